$wb = $excel.ActiveWorkbook

# --- Rename sheets: *-thanhpham -> *-linhkien -------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

$ws1.Name = "nhap-linhkien"
$ws2.Name = "xuat-linhkien"
$ws3.Name = "ton-linhkien"

# --- Sheet1 (nhap-linhkien): new header row, drop sample data row ----------
$ws1.Rows.Item(2).Delete()

$ws1.Range("A1").Value = "Part Number"
$ws1.Range("B1").Value = "Tên Hàng"
$ws1.Range("C1").Value = "Sổ Hợp Đồng"
$ws1.Range("D1").Value = "Sản Phẩm"
$ws1.Range("E1").Value = "Công Ty Nhập"
$ws1.Range("F1").Value = "Ngày Nhập"
$ws1.Range("G1").Value = "Đơn Vị Tính"
$ws1.Range("H1").Value = "Số Lượng"
$ws1.Range("I1").Value = "Đơn Giá"
$ws1.Range("J1").Value = "Thành Tiền"

$ws1.Range("A1:J1").ColumnWidth = 19.140625

# --- Sheet2 (xuat-linhkien): same new header row, no data rows -------------
$ws2.Range("A1").Value = "Part Number"
$ws2.Range("B1").Value = "Tên Hàng"
$ws2.Range("C1").Value = "Sổ Hợp Đồng"
$ws2.Range("D1").Value = "Sản Phẩm"
$ws2.Range("E1").Value = "Công Ty Nhập"
$ws2.Range("F1").Value = "Ngày Nhập"
$ws2.Range("G1").Value = "Đơn Vị Tính"
$ws2.Range("H1").Value = "Số Lượng"
$ws2.Range("I1").Value = "Đơn Giá"
$ws2.Range("J1").Value = "Thành Tiền"

$ws2.Range("A1:J1").ColumnWidth = 19.140625

# --- Sheet3 (ton-linhkien): new header row, drop sample data row -----------
$ws3.Rows.Item(2).Delete()

$ws3.Range("A1").Value = "Part Number"
$ws3.Range("B1").Value = "Số Lượng"
$ws3.Range("C1").Value = "Đơn Vị Tính"
